$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("1er Parcial")
$ws.Cells.Item(9, 5).Value = 36
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 100
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 8.6
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(10, 5).Value = 35
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 100
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 7.9
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(11, 5).Value = 25
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 100
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 8.4
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(12, 9).Value = 7.5
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(14, 9).Value = 7.4
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(16, 5).Value = 22
$ws.Cells.Item(16, 6).Value = 16
$ws.Cells.Item(16, 7).Value = 57.89
$ws.Cells.Item(16, 8).Value = 42.11
$ws.Cells.Item(16, 9).Value = 8
$ws.Cells.Item(16, 10).Value = 16
$ws.Cells.Item(16, 11).Value = 42.11
$ws.Cells.Item(17, 5).Value = 24
$ws.Cells.Item(17, 6).Value = 10
$ws.Cells.Item(17, 7).Value = 70.59
$ws.Cells.Item(17, 8).Value = 29.41
$ws.Cells.Item(17, 9).Value = 8.1
$ws.Cells.Item(17, 10).Value = 10
$ws.Cells.Item(17, 11).Value = 29.41

$ws = $wb.Worksheets.Item("2o Parcial")
$ws.Cells.Item(9, 5).Value = 36
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 100
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 8.5
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(10, 5).Value = 35
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 100
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 7.6
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(11, 5).Value = 25
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 100
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 8.300000000000001
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(12, 5).Value = 28
$ws.Cells.Item(12, 6).Value = 11
$ws.Cells.Item(12, 7).Value = 71.79000000000001
$ws.Cells.Item(12, 8).Value = 28.21
$ws.Cells.Item(12, 9).Value = 7.8
$ws.Cells.Item(12, 10).Value = 5
$ws.Cells.Item(12, 11).Value = 12.82
$ws.Cells.Item(13, 5).Value = 25
$ws.Cells.Item(13, 6).Value = 12
$ws.Cells.Item(13, 7).Value = 67.56999999999999
$ws.Cells.Item(13, 8).Value = 32.43
$ws.Cells.Item(13, 9).Value = 6.9
$ws.Cells.Item(13, 10).Value = 4
$ws.Cells.Item(13, 11).Value = 10.81
$ws.Cells.Item(14, 5).Value = 31
$ws.Cells.Item(14, 6).Value = 5
$ws.Cells.Item(14, 7).Value = 86.11
$ws.Cells.Item(14, 8).Value = 13.89
$ws.Cells.Item(14, 9).Value = 8.9
$ws.Cells.Item(14, 10).Value = 5
$ws.Cells.Item(14, 11).Value = 13.89
$ws.Cells.Item(15, 5).Value = 26
$ws.Cells.Item(15, 6).Value = 9
$ws.Cells.Item(15, 7).Value = 74.29000000000001
$ws.Cells.Item(15, 8).Value = 25.71
$ws.Cells.Item(15, 9).Value = 7.7
$ws.Cells.Item(15, 10).Value = 6
$ws.Cells.Item(15, 11).Value = 17.14
$ws.Cells.Item(16, 5).Value = 18
$ws.Cells.Item(16, 6).Value = 20
$ws.Cells.Item(16, 7).Value = 47.37
$ws.Cells.Item(16, 8).Value = 52.63
$ws.Cells.Item(16, 9).Value = 9.4
$ws.Cells.Item(16, 10).Value = 20
$ws.Cells.Item(16, 11).Value = 52.63
$ws.Cells.Item(17, 5).Value = 21
$ws.Cells.Item(17, 6).Value = 13
$ws.Cells.Item(17, 7).Value = 61.76
$ws.Cells.Item(17, 8).Value = 38.24
$ws.Cells.Item(17, 9).Value = 8.6
$ws.Cells.Item(17, 10).Value = 13
$ws.Cells.Item(17, 11).Value = 38.24

$ws = $wb.Worksheets.Item("3er Parcial")
$ws.Cells.Item(9, 5).Value = 36
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 100
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(10, 5).Value = 35
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 100
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 7.9
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(11, 5).Value = 25
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 100
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 0
$ws.Cells.Item(12, 5).Value = 32
$ws.Cells.Item(12, 6).Value = 7
$ws.Cells.Item(12, 7).Value = 82.05
$ws.Cells.Item(12, 8).Value = 17.95
$ws.Cells.Item(12, 9).Value = 7.7
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(13, 9).Value = 6.8
$ws.Cells.Item(14, 5).Value = 33
$ws.Cells.Item(14, 6).Value = 3
$ws.Cells.Item(14, 7).Value = 91.67
$ws.Cells.Item(14, 8).Value = 8.33
$ws.Cells.Item(14, 9).Value = 8.1
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(15, 9).Value = 7.5
$ws.Cells.Item(16, 5).Value = 22
$ws.Cells.Item(16, 6).Value = 16
$ws.Cells.Item(16, 7).Value = 57.89
$ws.Cells.Item(16, 8).Value = 42.11
$ws.Cells.Item(16, 9).Value = 8.5
$ws.Cells.Item(16, 10).Value = 16
$ws.Cells.Item(16, 11).Value = 42.11
$ws.Cells.Item(17, 5).Value = 24
$ws.Cells.Item(17, 6).Value = 10
$ws.Cells.Item(17, 7).Value = 70.59
$ws.Cells.Item(17, 8).Value = 29.41
$ws.Cells.Item(17, 9).Value = 8.5
$ws.Cells.Item(17, 10).Value = 10
$ws.Cells.Item(17, 11).Value = 29.41
